$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.758234
$ws.Range("H2").Value = 8.274702
$ws.Range("I2").Value = 0.1921318935545868
$ws.Range("J2").Value = 0.1921318935545868
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1148133333333333
$ws.Range("N2").Value = 0.34444
$ws.Range("O2").Value = 0.03343792635928704
$ws.Range("P2").Value = 0.03343792635928704
$ws.Range("Q2").Value = 0.3166820396533333
$ws.Range("R2").Value = 2.85013835688
$ws.Range("S2").Value = 0.006424492107948651
$ws.Range("T2").Value = 0.006424492107948651

$ws.Range("G3").Value = 2.758234
$ws.Range("H3").Value = 8.274702
$ws.Range("I3").Value = 0.1921318935545868
$ws.Range("J3").Value = 0.1921318935545868
$ws.Range("O3").Value = 0.9249645515654102
$ws.Range("P3").Value = 0.9249645515654102
$ws.Range("Q3").Value = 8.760102455199334
$ws.Range("R3").Value = 78.840922096794
$ws.Range("S3").Value = 0.1777151907631315
$ws.Range("T3").Value = 0.1777151907631315

$ws.Range("G4").Value = 2.758234
$ws.Range("H4").Value = 8.274702
$ws.Range("I4").Value = 0.1921318935545868
$ws.Range("J4").Value = 0.1921318935545868
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1428303333333333
$ws.Range("N4").Value = 0.428491
$ws.Range("O4").Value = 0.04159752207530271
$ws.Range("P4").Value = 0.04159752207530271
$ws.Range("Q4").Value = 0.3939594816313333
$ws.Range("R4").Value = 3.545635334682
$ws.Range("S4").Value = 0.007992210683506635
$ws.Range("T4").Value = 0.007992210683506635

$ws.Range("G5").Value = 7.720664
$ws.Range("I5").Value = 0.5378027367579149
$ws.Range("J5").Value = 0.5378027367579149
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1148133333333333
$ws.Range("N5").Value = 0.34444
$ws.Range("O5").Value = 0.03343792635928704
$ws.Range("P5").Value = 0.03343792635928704
$ws.Range("Q5").Value = 0.8864351693866667
$ws.Range("R5").Value = 7.977916524480001
$ws.Range("S5").Value = 0.01798300830753419
$ws.Range("T5").Value = 0.01798300830753419

$ws.Range("G6").Value = 7.720664
$ws.Range("I6").Value = 0.5378027367579149
$ws.Range("J6").Value = 0.5378027367579149
$ws.Range("O6").Value = 0.9249645515654102
$ws.Range("P6").Value = 0.9249645515654102
$ws.Range("R6").Value = 220.6862321904241
$ws.Range("S6").Value = 0.4974484672359351
$ws.Range("T6").Value = 0.4974484672359351

$ws.Range("G7").Value = 7.720664
$ws.Range("I7").Value = 0.5378027367579149
$ws.Range("J7").Value = 0.5378027367579149
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.1428303333333333
$ws.Range("N7").Value = 0.428491
$ws.Range("O7").Value = 0.04159752207530271
$ws.Range("P7").Value = 0.04159752207530271
$ws.Range("Q7").Value = 1.102745012674667
$ws.Range("R7").Value = 9.924705114072001
$ws.Range("S7").Value = 0.02237126121444558
$ws.Range("T7").Value = 0.02237126121444558

$ws.Range("G8").Value = 3.877042333333333
$ws.Range("H8").Value = 11.631127
$ws.Range("I8").Value = 0.2700653696874982
$ws.Range("J8").Value = 0.2700653696874982
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1148133333333333
$ws.Range("N8").Value = 0.34444
$ws.Range("O8").Value = 0.03343792635928704
$ws.Range("P8").Value = 0.03343792635928704
$ws.Range("Q8").Value = 0.4451361537644444
$ws.Range("R8").Value = 4.00622538388
$ws.Range("S8").Value = 0.009030425943804198
$ws.Range("T8").Value = 0.009030425943804198

$ws.Range("G9").Value = 3.877042333333333
$ws.Range("H9").Value = 11.631127
$ws.Range("I9").Value = 0.2700653696874982
$ws.Range("J9").Value = 0.2700653696874982
$ws.Range("O9").Value = 0.9249645515654102
$ws.Range("P9").Value = 0.9249645515654102
$ws.Range("Q9").Value = 12.31341795625211
$ws.Range("R9").Value = 110.820761606269
$ws.Range("S9").Value = 0.2498008935663435
$ws.Range("T9").Value = 0.2498008935663435

$ws.Range("G10").Value = 3.877042333333333
$ws.Range("H10").Value = 11.631127
$ws.Range("I10").Value = 0.2700653696874982
$ws.Range("J10").Value = 0.2700653696874982
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.1428303333333333
$ws.Range("N10").Value = 0.428491
$ws.Range("O10").Value = 0.04159752207530271
$ws.Range("P10").Value = 0.04159752207530271
$ws.Range("Q10").Value = 0.5537592488174444
$ws.Range("R10").Value = 4.983833239357
$ws.Range("S10").Value = 0.01123405017735049
$ws.Range("T10").Value = 0.01123405017735049
